$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 236, shifting existing rows 236-279 down to 237-280.
$ws.Rows.Item(236).Insert()

# Populate the newly inserted row 236 with the new weekly price record.
$ws.Cells.Item(236, 1).Value = 3
$ws.Cells.Item(236, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(236, 3).Value = "Coquimbo"
$ws.Cells.Item(236, 4).Value = 44522
$ws.Cells.Item(236, 5).Value = 5
$ws.Cells.Item(236, 6).Value = 100112028
$ws.Cells.Item(236, 7).Value = "Sandia"
$ws.Cells.Item(236, 8).Value = "Sin especificar"
$ws.Cells.Item(236, 9).Value = "Primera"
$ws.Cells.Item(236, 10).Value = 160
$ws.Cells.Item(236, 11).Value = 700
$ws.Cells.Item(236, 12).Value = 700
$ws.Cells.Item(236, 13).Value = 700
$ws.Cells.Item(236, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(236, 15).Value = "Perú"
$ws.Cells.Item(236, 16).Value = 700
$ws.Cells.Item(236, 17).Value = 1
$ws.Cells.Item(236, 18).Value = "Hortaliza"
